$d = $word.ActiveDocument

# The last paragraph of the document currently holds the diary date line
# "2023年3月22日星期三 小雨" (with the _GoBack bookmark attached to it).
# The edit:
#   1. Moves that date text into a brand-new paragraph inserted just above
#      (formatted like the other diary-date paragraphs: rFonts hint="eastAsia").
#   2. Leaves the original paragraph in place (same pPr / bookmark) but
#      replaces its run text with the new diary entry "使用Git简单又便捷。".

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range

# Exclude the trailing paragraph mark from the range so we only replace the
# run text that precedes it (keeps the paragraph mark - and the bookmark
# anchored at its position - intact for re-use below).
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)

$oldDateText = $textRange.Text
$newEntryText = "使用Git简单又便捷。"

# Defensively XML-escape the text we splice into the raw OOXML fragment.
function Escape-Xml([string]$s) {
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}
$oldDateTextXml = Escape-Xml $oldDateText
$newEntryTextXml = Escape-Xml $newEntryText

$newParagraphXml = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>' + $oldDateTextXml + '</w:t></w:r></w:p>'

$updatedOriginalParagraphXml = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>' + $newEntryTextXml + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$bodyXml = $newParagraphXml + $updatedOriginalParagraphXml

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$textRange.InsertXML($packageXml)
